$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6 and Row 7: fill in D/G/H columns (Automatable="A", Steps, Expected Behaviour) ---
# Order matters: new shared strings are appended in the order the distinct values are
# first written, and the target file expects (in append order) H7, G7, G6, H6.

$h7 = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Native Tabbar JS Test
};
validate3
{
validate_Text_Exists=VT200-0576
};
validate4
{
validate_Screenshot=VT200_0576_new
};
validate5
{
validate_Text_Exists=VT200-0596
};
validate6
{
validate_Screenshot=VT200_0596
};
"@

$g7 = @"
wait(5);
validate1;
link_Click(tabbar_test_link);
validate2;
SelectTestToRun(VT200_0576_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(5);
TakeScreenshot(VT200_0576_new);
validate4;
SwitchApp(NATIVE_APP);
wait(2);
link_Click(VT200_0576_mainpage_xpath);
wait(2);
SwitchApp(WEBVIEW);
SelectTestToRun(VT200_0596_string);
ClickRunTest(runtest_top_xpath);
validate5;
ClickRunTest(runtest_bottom_xpath);
wait(5);
TakeNativeScreenshot(VT200_0596);
validate6;
"@

$g6 = @"
wait(5);
validate1;
link_Click(tabbar_test_link);
validate2;
SelectTestToRun(VT200_0576_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(5);
TakeScreenshot(VT200_0576_new);
validate4;
SwitchApp(NATIVE_APP);
wait(2);
link_Click(VT200_0576_mainpage_xpath);
wait(2);
SwitchApp(WEBVIEW);
SelectTestToRun(VT200_0595_string);
ClickRunTest(runtest_top_xpath);
validate5;
ClickRunTest(runtest_bottom_xpath);
wait(5);
TakeNativeScreenshot(VT200_0595);
validate6;

"@

$h6 = @"
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Native Tabbar JS Test
};
validate3
{
validate_Text_Exists=VT200-0576
};
validate4
{
validate_Screenshot=VT200_0576_new
};
validate5
{
validate_Text_Exists=VT200-0595
};
validate6
{
validate_Screenshot=VT200_0595
};

"@

$ws.Range("D7").Value = "A"
$ws.Range("H7").Value = $h7
$ws.Range("G7").Value = $g7

$ws.Range("D6").Value = "A"
$ws.Range("G6").Value = $g6
$ws.Range("H6").Value = $h6

# --- Row heights ---
$ws.Rows.Item(6).RowHeight = 282
$ws.Rows.Item(7).RowHeight = 270.75

# --- Column G width ---
$ws.Columns.Item(7).ColumnWidth = 35.6667

# --- New row 14 / D14: empty cell carrying a border-only style ---
$ws.Range("I1").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Borders.Item(9).LineStyle = -4142
$ws.Range("D14").ClearContents()

# --- Selection moves to D1 ---
$ws.Range("D1").Select()

Write-Output "done"
